# Append Chris Gayle IPL match rows 9-15 to sheet1 (Kings XI Punjab data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = ' Abu Dhabi'
$ws.Range("B9").Value = ' October 30 2020'
$ws.Range("C9").Value = 'Royals won by 7 wickets (with 15 balls remaining)'
$ws.Range("D9").Value = 'Kings XI Punjab'
$ws.Range("E9").Value = 'Rajasthan Royals'
$ws.Range("F9").Value = 'Chris Gayle '
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '99'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '63'
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = '6'
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = '8'
$ws.Range("K9").NumberFormat = "@"
$ws.Range("K9").Value = '157.14'

# Row 10
$ws.Range("A10").Value = ' Sharjah'
$ws.Range("B10").Value = ' October 26 2020'
$ws.Range("C10").Value = 'Kings XI won by 8 wickets (with 7 balls remaining)'
$ws.Range("D10").Value = 'Kings XI Punjab'
$ws.Range("E10").Value = 'Kolkata Knight Riders'
$ws.Range("F10").Value = 'Chris Gayle '
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '51'
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = '29'
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = '2'
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = '5'
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = '175.86'

# Row 11
$ws.Range("A11").Value = ' Abu Dhabi'
$ws.Range("B11").Value = ' November 01 2020'
$ws.Range("C11").Value = 'Super Kings won by 9 wickets (with 7 balls remaining)'
$ws.Range("D11").Value = 'Kings XI Punjab'
$ws.Range("E11").Value = 'Chennai Super Kings'
$ws.Range("F11").Value = 'Chris Gayle '
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '12'
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '19'
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = '0'
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = '0'
$ws.Range("K11").NumberFormat = "@"
$ws.Range("K11").Value = '63.15'

# Row 12
$ws.Range("A12").Value = ' Sharjah'
$ws.Range("B12").Value = ' October 15 2020'
$ws.Range("C12").Value = 'Kings XI won by 8 wickets'
$ws.Range("D12").Value = 'Kings XI Punjab'
$ws.Range("E12").Value = 'Royal Challengers Bangalore'
$ws.Range("F12").Value = 'Chris Gayle '
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '53'
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = '45'
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = '1'
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = '5'
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = '117.77'

# Row 13
$ws.Range("A13").Value = ' Dubai (DSC)'
$ws.Range("B13").Value = ' October 20 2020'
$ws.Range("C13").Value = 'Kings XI won by 5 wickets (with 6 balls remaining)'
$ws.Range("D13").Value = 'Kings XI Punjab'
$ws.Range("E13").Value = 'Delhi Capitals'
$ws.Range("F13").Value = 'Chris Gayle '
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '29'
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = '13'
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = '3'
$ws.Range("J13").NumberFormat = "@"
$ws.Range("J13").Value = '2'
$ws.Range("K13").NumberFormat = "@"
$ws.Range("K13").Value = '223.07'

# Row 14
$ws.Range("A14").Value = ' Dubai (DSC)'
$ws.Range("B14").Value = ' October 18 2020'
$ws.Range("C14").Value = 'Match tied (Kings XI won the one-over eliminator)'
$ws.Range("D14").Value = 'Kings XI Punjab'
$ws.Range("E14").Value = 'Mumbai Indians'
$ws.Range("F14").Value = 'Chris Gayle '
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '24'
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = '21'
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = '1'
$ws.Range("J14").NumberFormat = "@"
$ws.Range("J14").Value = '2'
$ws.Range("K14").NumberFormat = "@"
$ws.Range("K14").Value = '114.28'

# Row 15
$ws.Range("A15").Value = ' Dubai (DSC)'
$ws.Range("B15").Value = ' October 24 2020'
$ws.Range("C15").Value = 'Kings XI won by 12 runs'
$ws.Range("D15").Value = 'Kings XI Punjab'
$ws.Range("E15").Value = 'Sunrisers Hyderabad'
$ws.Range("F15").Value = 'Chris Gayle '
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '20'
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = '20'
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = '2'
$ws.Range("J15").NumberFormat = "@"
$ws.Range("J15").Value = '1'
$ws.Range("K15").NumberFormat = "@"
$ws.Range("K15").Value = '100.00'
